# edit.ps1 - applies the "Abbreviations" slide split edit described by the diff:
#   * Duplicates the "Abbreviations" slide (originally slide index 3) and moves the
#     duplicate in front of the original, producing a new slide (PowerPoint assigns
#     it SlideID 278, matching the target sldIdLst) right before the original.
#   * The new (duplicate, now at position 3) slide becomes "Abbreviations - business":
#       - keeps only the business-terms column, with its text expanded / reworded
#       - the two "technology" columns are removed
#       - the remaining column is widened slightly
#   * The original slide (now at position 4) becomes "Abbreviations - technology stack":
#       - the business-terms column is removed
#       - the two technology columns are re-positioned to fill the freed space

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: convert EMU (OOXML units, 914400 per inch) to points (PowerPoint's
# COM Left/Top/Width/Height unit, 72 per inch), nudged so that PowerPoint's
# internal float rounding lands back on the exact EMU value on save.
# ---------------------------------------------------------------------------
function EMU([double]$emu) {
    return ($emu / 12700.0) + 0.00003
}

function Get-ShapeById($slide, [int]$id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Locate the existing "Abbreviations" slide (has title text exactly
# "Abbreviations" and three "Content Placeholder 2" shapes).
# ---------------------------------------------------------------------------
$abbrevIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($sh in $slide.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "Abbreviations") {
                $abbrevIndex = $i
                break
            }
        }
    }
    if ($abbrevIndex -gt 0) { break }
}

$original = $p.Slides.Item($abbrevIndex)

# Duplicate it and move the duplicate to sit right before the original.
$dup = $original.Duplicate()
$dup.MoveTo($abbrevIndex)

$businessSlide = $p.Slides.Item($abbrevIndex)
$techSlide = $p.Slides.Item($abbrevIndex + 1)

# ---------------------------------------------------------------------------
# New slide -> "Abbreviations - business"
# ---------------------------------------------------------------------------
foreach ($sh in $businessSlide.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "Abbreviations") {
        $sh.TextFrame.TextRange.Text = "Abbreviations - business"
    }
}

# Remove the two "technology" columns (ids 4 and 5 on the original slide).
$biz_tech1 = Get-ShapeById $businessSlide 4
if ($biz_tech1 -ne $null) { $biz_tech1.Delete() }
$biz_tech2 = Get-ShapeById $businessSlide 5
if ($biz_tech2 -ne $null) { $biz_tech2.Delete() }

# Widen the remaining (business-terms) column and rewrite its content.
$bizCol = Get-ShapeById $businessSlide 3
$bizCol.Left = EMU 296586
$bizCol.Top = EMU 914400
$bizCol.Width = EMU 3969026
$bizCol.Height = EMU 4724400

$bizTextRange = $bizCol.TextFrame.TextRange
$bizTextRange.Text = "QMS - Quality Management System `rQC - Quality Control`rQA - Quality Assurance`rNCR - Non-Conformance Request / Report`rPOC  - Proof Of Concept`rRFI - Request For Information`rCAR - Corrective Action Request`rPAR - Preventive Action Request`rOBS"
$bizTextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Original slide (now pushed one position later) -> "Abbreviations - technology stack"
# ---------------------------------------------------------------------------
foreach ($sh in $techSlide.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "Abbreviations") {
        $sh.TextFrame.TextRange.Text = "Abbreviations - technology stack"
    }
}

# Remove the business-terms column (id 3) from the technology-stack slide.
$techBiz = Get-ShapeById $techSlide 3
if ($techBiz -ne $null) { $techBiz.Delete() }

# Re-position the two technology columns to fill the freed space.
$techCol1 = Get-ShapeById $techSlide 4
$techCol1.Left = EMU 1827212
$techCol1.Top = EMU 914400

$techCol2 = Get-ShapeById $techSlide 5
$techCol2.Left = EMU 5256212
$techCol2.Top = EMU 849464

Write-Host "Done. Slides.Count =" $p.Slides.Count
